$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.255432847381883
$ws.Range("C2").Value = 0.2318737208821346
$ws.Range("D2").Value = 0.01861625622255048
$ws.Range("F2").Value = 0.8596978838786953
$ws.Range("G2").Value = 0.7134475634870228
$ws.Range("H2").Value = 0.7629230913686342
$ws.Range("I2").Value = 0.6978593141990928
$ws.Range("L2").Value = 0.3010685567896303
$ws.Range("N2").Value = 1.114843355685437
$ws.Range("B3").Value = 1.138903253550609
$ws.Range("C3").Value = 0.2121045258404592
$ws.Range("D3").Value = 0.01838677076083073
$ws.Range("F3").Value = 0.8443165825569849
$ws.Range("G3").Value = 0.6981026244201018
$ws.Range("H3").Value = 0.7611619914334682
$ws.Range("I3").Value = 0.6992619156683588
$ws.Range("L3").Value = 0.2896779345044536
$ws.Range("N3").Value = 1.129490372969432
$ws.Range("B4").Value = 1.067653821679755
$ws.Range("C4").Value = 0.1998722072248142
$ws.Range("D4").Value = 0.01824699074694536
$ws.Range("F4").Value = 0.835534889404812
$ws.Range("G4").Value = 0.6892925885629353
$ws.Range("H4").Value = 0.7605618622210244
$ws.Range("I4").Value = 0.7006216511125416
$ws.Range("L4").Value = 0.2828630492595465
$ws.Range("N4").Value = 1.138982913117953
$ws.Range("B5").Value = 1.03869535343938
$ws.Range("C5").Value = 0.1948638461757639
$ws.Range("D5").Value = 0.0181903196149733
$ws.Range("F5").Value = 0.8321223300683585
$ws.Range("G5").Value = 0.685855581282155
$ws.Range("H5").Value = 0.7604381048565472
$ws.Range("I5").Value = 0.7013007945368983
$ws.Range("L5").Value = 0.2801308851404372
$ws.Range("N5").Value = 1.142976647556381
$ws.Range("B6").Value = 1.033891450160922
$ws.Range("C6").Value = 0.1940307877740395
$ws.Range("D6").Value = 0.0181809271759974
$ws.Range("F6").Value = 0.8315656906901268
$ws.Range("G6").Value = 0.6852940987292584
$ws.Range("H6").Value = 0.7604248452672806
$ws.Range("I6").Value = 0.7014211085701731
$ws.Range("L6").Value = 0.279679925099245
$ws.Range("N6").Value = 1.143647378464173
$ws.Range("B7").Value = 1.067262967883437
$ws.Range("C7").Value = 0.1998047581916467
$ws.Range("D7").Value = 0.01824622527409758
$ws.Range("F7").Value = 0.8354881948105799
$ws.Range("G7").Value = 0.6892456166856959
$ws.Range("H7").Value = 0.7605597043389025
$ws.Range("I7").Value = 0.7006303044221838
$ws.Range("L7").Value = 0.2828260203797157
$ws.Range("N7").Value = 1.139036266198822
$ws.Range("B8").Value = 1.215191553143768
$ws.Range("C8").Value = 0.2250768165724537
$ws.Range("D8").Value = 0.01853690049803802
$ws.Range("F8").Value = 0.8542565814577756
$ws.Range("G8").Value = 0.7080292213697845
$ws.Range("H8").Value = 0.7622158689969467
$ws.Range("I8").Value = 0.6982392662588879
$ws.Range("L8").Value = 0.2971038806215773
$ws.Range("N8").Value = 1.119789985597755
$ws.Range("B9").Value = 1.507643121222429
$ws.Range("C9").Value = 0.2738905910703693
$ws.Range("D9").Value = 0.01911554982968511
$ws.Range("F9").Value = 0.8963452981688533
$ws.Range("G9").Value = 0.7497525581070477
$ws.Range("H9").Value = 0.7692920135211097
$ws.Range("I9").Value = 0.6975213846666861
$ws.Range("L9").Value = 0.3265273782536298
$ws.Range("N9").Value = 1.086012522016407
$ws.Range("B10").Value = 1.723947375444482
$ws.Range("C10").Value = 0.309304591979668
$ws.Range("D10").Value = 0.01954561440838276
$ws.Range("F10").Value = 0.9305317590895612
$ws.Range("G10").Value = 0.78343869186088
$ws.Range("H10").Value = 0.7768410416688596
$ws.Range("I10").Value = 0.6994372011537919
$ws.Range("L10").Value = 0.3490220623482543
$ws.Range("N10").Value = 1.063617232799331
$ws.Range("B11").Value = 1.82266382447159
$ws.Range("C11").Value = 0.3253189114414283
$ws.Range("D11").Value = 0.0197422680416679
$ws.Range("F11").Value = 0.9468022762999055
$ws.Range("G11").Value = 0.7994331653758309
$ws.Range("H11").Value = 0.7807892714290006
$ws.Range("I11").Value = 0.7008442238030881
$ws.Range("L11").Value = 0.3594480197557317
$ws.Range("N11").Value = 1.053955718464763
$ws.Range("B12").Value = 1.860090587432637
$ws.Range("C12").Value = 0.3313693774828721
$ws.Range("D12").Value = 0.01981687549299238
$ws.Range("F12").Value = 0.9530675533901274
$ws.Range("G12").Value = 0.8055871077018253
$ws.Range("H12").Value = 0.7823585605832477
$ws.Range("I12").Value = 0.7014544142489072
$ws.Range("L12").Value = 0.3634239238421486
$ws.Range("N12").Value = 1.050372941699361
$ws.Range("B13").Value = 1.852028076597776
$ws.Range("C13").Value = 0.3300669163089651
$ws.Range("D13").Value = 0.01980080137063212
$ws.Range("F13").Value = 0.9517135794359746
$ws.Range("G13").Value = 0.8042574116523724
$ws.Range("H13").Value = 0.782017282933765
$ws.Range("I13").Value = 0.701319551248119
$ws.Range("L13").Value = 0.3625664033941689
$ws.Range("N13").Value = 1.0511411816546
$ws.Range("B14").Value = 1.82574204987236
$ws.Range("C14").Value = 0.3258169639603068
$ws.Range("D14").Value = 0.01974840329534189
$ws.Range("F14").Value = 0.9473156369232925
$ws.Range("G14").Value = 0.7999375018301578
$ws.Range("H14").Value = 0.780916889684562
$ws.Range("I14").Value = 0.7008928715377749
$ws.Range("L14").Value = 0.3597745618398847
$ws.Range("N14").Value = 1.053659440601091
$ws.Range("B15").Value = 1.809646933790987
$ws.Range("C15").Value = 0.3232119456177713
$ws.Range("D15").Value = 0.01971632585691552
$ws.Range("F15").Value = 0.9446353300964034
$ws.Range("G15").Value = 0.7973041111629584
$ws.Range("H15").Value = 0.780252534885733
$ws.Range("I15").Value = 0.7006416064406764
$ws.Range("L15").Value = 0.3580681032537569
$ws.Range("N15").Value = 1.055211826940869
$ws.Range("B16").Value = 1.717502363057974
$ws.Range("C16").Value = 0.3082560888176431
$ws.Range("D16").Value = 0.01953278245686718
$ws.Range("F16").Value = 0.929482960780021
$ws.Range("G16").Value = 0.7824069728010556
$ws.Range("H16").Value = 0.776593382486169
$ws.Range("I16").Value = 0.6993560568186012
$ws.Range("L16").Value = 0.3483445909800054
$ws.Range("N16").Value = 1.064259237817293
$ws.Range("B17").Value = 1.661055656027258
$ws.Range("C17").Value = 0.2990566117811397
$ws.Range("D17").Value = 0.01942043949134842
$ws.Range("F17").Value = 0.9203720365046877
$ws.Range("G17").Value = 0.7734403331351416
$ws.Range("H17").Value = 0.7744804699469228
$ws.Range("I17").Value = 0.6987048373259839
$ws.Range("L17").Value = 0.3424290124747813
$ws.Range("N17").Value = 1.069944419645719
$ws.Range("B18").Value = 1.62861903592966
$ws.Range("C18").Value = 0.2937563122814595
$ws.Range("D18").Value = 0.01935591874021725
$ws.Range("F18").Value = 0.9151993061228438
$ws.Range("G18").Value = 0.7683460220845859
$ws.Range("H18").Value = 0.7733135567890201
$ws.Range("I18").Value = 0.6983806550595517
$ws.Range("L18").Value = 0.3390446914973637
$ws.Range("N18").Value = 1.073263905773718
$ws.Range("B19").Value = 1.6176417370786
$ws.Range("C19").Value = 0.2919601770823874
$ws.Range("D19").Value = 0.01933408981833651
$ws.Range("F19").Value = 0.9134595060171762
$ws.Range("G19").Value = 0.7666319836584989
$ws.Range("H19").Value = 0.7729267615473532
$ws.Range("I19").Value = 0.6982795342370736
$ws.Range("L19").Value = 0.3379019361581754
$ws.Range("N19").Value = 1.074396327058672
$ws.Range("B20").Value = 1.667061402054401
$ws.Range("C20").Value = 0.3000368440224292
$ws.Range("D20").Value = 0.01943238871244546
$ws.Range("F20").Value = 0.9213349058825031
$ws.Range("G20").Value = 0.7743883149022679
$ws.Range("H20").Value = 0.7747003843339542
$ws.Range("I20").Value = 0.698768943685927
$ws.Range("L20").Value = 0.3430568556025833
$ws.Range("N20").Value = 1.069334096064001
$ws.Range("B21").Value = 1.833461678698654
$ws.Range("C21").Value = 0.3270656531522036
$ws.Range("D21").Value = 0.01976379016284469
$ws.Range("F21").Value = 0.9486045922828623
$ws.Range("G21").Value = 0.8012037204300384
$ws.Range("H21").Value = 0.7812380866417641
$ws.Range("I21").Value = 0.7010160945940456
$ws.Range("L21").Value = 0.3605938378199767
$ws.Range("N21").Value = 1.052917707623727
$ws.Range("B22").Value = 1.942475717588366
$ws.Range("C22").Value = 0.344650029672863
$ws.Range("D22").Value = 0.01998118719067676
$ws.Range("F22").Value = 0.9670332972267346
$ws.Range("G22").Value = 0.81929594905165
$ws.Range("H22").Value = 0.7859433485556622
$ws.Range("I22").Value = 0.7029359519661398
$ws.Range("L22").Value = 0.3722174584536759
$ws.Range("N22").Value = 1.042630716663215
$ws.Range("B23").Value = 1.8842691843343
$ws.Range("C23").Value = 0.3352722999883326
$ws.Range("D23").Value = 0.01986508672276344
$ws.Range("F23").Value = 0.9571418655432637
$ws.Range("G23").Value = 0.8095876686156487
$ws.Range("H23").Value = 0.7833924060057882
$ws.Range("I23").Value = 0.7018698771219647
$ws.Range("L23").Value = 0.3659988512168866
$ws.Range("N23").Value = 1.048080573783825
$ws.Range("B24").Value = 1.664346156422425
$ws.Range("C24").Value = 0.2995937165299551
$ws.Range("D24").Value = 0.01942698626154638
$ws.Range("F24").Value = 0.9208993893937816
$ws.Range("G24").Value = 0.7739595432334454
$ws.Range("H24").Value = 0.7746008120679448
$ws.Range("I24").Value = 0.6987398048141529
$ws.Range("L24").Value = 0.3427729559518014
$ws.Range("N24").Value = 1.069609864497714
$ws.Range("B25").Value = 1.428274284682516
$ws.Range("C25").Value = 0.260764432452504
$ws.Range("D25").Value = 0.01895812049994205
$ws.Range("F25").Value = 0.884389101704059
$ws.Range("G25").Value = 0.7379365817666894
$ws.Range("H25").Value = 0.7669661955138594
$ws.Range("I25").Value = 0.6972881597036462
$ws.Range("L25").Value = 0.3184141245500882
$ws.Range("N25").Value = 1.094725390265339
